# Add a "Save" column (column H) to the s_vals worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header (copy formatting from the neighboring "sum" header so the new
# column matches the existing bold/bordered/centered header style)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for the new "Save" column (row -> value)
$saveValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
